# VVencidos.xlsx — "Add files via upload"
#
# The refreshed export (re-uploaded workbook) fills in the running
# "Entidade #" counter in column A for rows 782-877 (it was already present
# for rows 2-781, i.e. row - 2), and leaves the grid scrolled/selected near
# the bottom of the sheet (selection on C884) instead of the original
# top-of-sheet selection (G2:G877).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A ("Id") was missing for the last 96 data rows (782-877) of the
# sheet; back-fill it following the existing pattern used by every row
# above it: A<row> = row - 2.
for ($r = 782; $r -le 877; $r++) {
    $ws.Cells.Item($r, 1).Value = $r - 2
}

# Leave the workbook with the selection where the author left it after
# uploading/reviewing the new rows — near the bottom of the sheet.
$ws.Range("C884").Select()
